$d = $word.ActiveDocument

# --- 1. "List others who contribute..." -- drop the gramStart/gramEnd
#        proofing markers around "List" (no visible text change). We
#        re-find/replace the phrase (including the leading space) so the
#        runs collapse back into one another exactly like an accepted
#        proofing pass in real Word would.
$d.Content.Find.Execute(
    " List others who contribute to a section as Helpers.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " List others who contribute to a section as Helpers.", 2) | Out-Null

# --- 2. "...is your client, the specific person..." -> "...is the specific
#        person...", and drop the trailing space that followed the sentence.
$d.Content.Find.Execute(
    "is your client, the specific person who can decide whether to implement your recommendations. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "is the specific person who can decide whether to implement your recommendations.", 2) | Out-Null

# --- 3. Remove the extra "before" spacing on the Heading2 paragraphs that
#        had an explicit 6pt (120 twips) override -- they fall back to the
#        style's own 12pt (240 twips) "before" spacing.
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Heading 2" -and $p.Format.SpaceBefore -eq 6) {
        $p.Format.SpaceBefore = 12
    }
}

# --- 4. "(You'll address the Recommendation Report to this person)" ->
#        "(You'll address the report to this person)"
$d.Content.Find.Execute(
    "(You’ll address the Recommendation Report to this person)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(You’ll address the report to this person)", 2) | Out-Null

# --- 5. " for your Recommendation Report? " -> " for your recommendation report? "
$d.Content.Find.Execute(
    " for your Recommendation Report? ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " for your recommendation report? ", 2) | Out-Null

# --- 6. Remove the whole "TikTok:" list paragraph.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "TikTok:`r") {
        $p.Range.Delete()
        break
    }
}

# --- 7. "generally to outside consultants and/or" -> "generally to students and/or"
$d.Content.Find.Execute(
    "generally to outside consultants and/or ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "generally to students and/or ", 2) | Out-Null

# --- 8. Drop the "Why do they care about it?" sentence.
$d.Content.Find.Execute(
    " (e.g., How much do they care? What do they care about? Why do they care about it?)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " (e.g., How much do they care? What do they care about?)", 2) | Out-Null

# --- 9. "...get them involved and invested in your recommendations?" -- drop
#        the gramStart/gramEnd proofing markers around "invested" (no
#        visible text change).
$d.Content.Find.Execute(
    " you use to get them involved and invested in your recommendations?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " you use to get them involved and invested in your recommendations?", 2) | Out-Null

# --- 10. Trim the run of 17 empty trailing paragraphs plus the final
#         tab-stop-only paragraph right before the section break, leaving
#         the last "real" (empty Boxed-style) paragraph intact.
$total = $d.Paragraphs.Count
$lastKeep = $null
for ($i = $total; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -ne "`r" -and $p.Range.Text -ne "`t`r") {
        $lastKeep = $i
        break
    }
}
if ($lastKeep -ne $null -and $lastKeep -lt $total) {
    $startPara = $d.Paragraphs.Item($lastKeep + 1)
    $endPara = $d.Paragraphs.Item($total)
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}

# --- 11. Footer: "Copyright © 2019–2025 Traci Gardner" -> "...2023..."
$footer = $d.Sections(1).Footers(1)
$footer.Range.Find.Execute(
    "Copyright © 2019–2025 Traci Gardner",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Copyright © 2019–2023 Traci Gardner", 2) | Out-Null

# --- 12. Footer: collapse the CC-license run fragments (proofing marks
#         around "NonCommercial"/"ShareAlike") back into one run.
$footer.Range.Find.Execute(
    "CC Attribution-NonCommercial-ShareAlike 4.0 International",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "CC Attribution-NonCommercial-ShareAlike 4.0 International", 2) | Out-Null
